# Sample Project / Main.xlsx - "Rules" sheet, row 11 (the R40 / 22-23 rule
# row) had its Rule-name cell (B11) changed from the label "R40" to the
# literal text "1".
#
# A leading apostrophe forces Excel to store the value as text (shared
# string) instead of re-interpreting the digit "1" as a number, matching
# the workbook's <c t="s"> cell type for B11 in the saved file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = "'1"
